# Modify categories for section 19 (sheets "19.1_L2_Cases_Apr-Sep2025" and
# "19.2_L3_Cases_Apr-Sep2025"): break the combined "Portals" and
# "Public Safety" groupings into their individual platform categories, and
# add several new, previously-unlisted categories, recomputing the Total row.

$wb = $excel.ActiveWorkbook

$sheetNames = @("19.1_L2_Cases_Apr-Sep2025", "19.2_L3_Cases_Apr-Sep2025")

# Row data per sheet: Platform Group, Case Count, % of Total
$data19_1 = @(
    @("MAC+", 1708, "50.1%"),
    @("TAP", 407, "11.9%"),
    @("MGI", 297, "8.7%"),
    @("GIFR", 132, "3.9%"),
    @("USB", 128, "3.8%"),
    @("GEARS", 107, "3.1%"),
    @("LMS", 47, "1.4%"),
    @("FAS", 40, "1.2%"),
    @("CORE PATHWAY", 6, "0.2%"),
    @("RLH Online", 1, "0.0%"),
    @("Online Storefront (Shopify)", 1, "0.0%"),
    @("API Integration (Janus)", 0, "0.0%"),
    @("Total", 2874, "84.3%")
)

$data19_2 = @(
    @("MAC+", 203, "6.0%"),
    @("TAP", 104, "3.0%"),
    @("MGI", 114, "3.3%"),
    @("GIFR", 0, "0.0%"),
    @("USB", 6, "0.2%"),
    @("GEARS", 1, "0.0%"),
    @("LMS", 2, "0.1%"),
    @("FAS", 2, "0.1%"),
    @("CORE PATHWAY", 4, "0.1%"),
    @("RLH Online", 0, "0.0%"),
    @("Online Storefront (Shopify)", 0, "0.0%"),
    @("API Integration (Janus)", 0, "0.0%"),
    @("Total", 436, "12.8%")
)

$dataBySheet = @{
    "19.1_L2_Cases_Apr-Sep2025" = $data19_1
    "19.2_L3_Cases_Apr-Sep2025" = $data19_2
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $dataBySheet[$sheetName]

    # Data now spans rows 2..14 (13 rows) instead of the original 2..5.
    $lastRow = 1 + $rows.Count

    # Ensure column C keeps being stored as plain text (e.g. "50.1%") rather
    # than being auto-converted into a numeric percentage by Excel.
    $ws.Range("C2:C$lastRow").NumberFormat = "@"

    $r = 2
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $r = $r + 1
    }
}
